# Updates the cryptocurrency price/volume table to reflect the refreshed
# market data (and the BitcoinSV/EnergySwap row swap) captured by the
# scheduled GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text such as '30.552.06' or '0.000007587' that would
# be silently reinterpreted as a number (and lose its exact formatting) if
# assigned as-is, so force the column to Text format first.
$ws.Range('D2:D51').NumberFormat = '@'

# --- Coin / Link swap (rows 48 and 49 traded places) ---
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'

# --- Price (column D) ---
$ws.Range('D2').Value = '30.552.06'
$ws.Range('D3').Value = '1.915.18'
$ws.Range('D5').Value = '244.13'
$ws.Range('D7').Value = '0.4844'
$ws.Range('D8').Value = '0.2900'
$ws.Range('D9').Value = '0.06863'
$ws.Range('D10').Value = '111.17'
$ws.Range('D11').Value = '19.31'
$ws.Range('D12').Value = '1.920.35'
$ws.Range('D13').Value = '0.07568'
$ws.Range('D14').Value = '5.343'
$ws.Range('D15').Value = '0.6710'
$ws.Range('D16').Value = '293.42'
$ws.Range('D17').Value = '30.542.44'
$ws.Range('D20').Value = '0.000007587'
$ws.Range('D21').Value = '2.166.17'
$ws.Range('D22').Value = '5.505'
$ws.Range('D25').Value = '9.454'
$ws.Range('D26').Value = '164.85'
$ws.Range('D27').Value = '20.24'
$ws.Range('D28').Value = '2.094'
$ws.Range('D30').Value = '1.433'
$ws.Range('D31').Value = '4.137'
$ws.Range('D32').Value = '4.056'
$ws.Range('D33').Value = '0.04984'
$ws.Range('D34').Value = '0.7360'
$ws.Range('D35').Value = '1.135'
$ws.Range('D36').Value = '1.000'
$ws.Range('D37').Value = '2.709'
$ws.Range('D38').Value = '0.02027'
$ws.Range('D39').Value = '2.683'
$ws.Range('D40').Value = '2.016'
$ws.Range('D41').Value = '109.56'
$ws.Range('D42').Value = '0.4436'
$ws.Range('D43').Value = '0.8622'
$ws.Range('D44').Value = '5.834'
$ws.Range('D46').Value = '69.45'
$ws.Range('D47').Value = '7.205'
$ws.Range('D48').Value = '9.240'
$ws.Range('D49').Value = '48.14'
$ws.Range('D50').Value = '0.2548'
$ws.Range('D51').Value = '0.1226'

# --- Volume(1h) (column E) ---
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').Value = '  +2.17%  '
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('E9').Value = '  +0.09%  '
$ws.Range('E10').Value = '  +5.62%  '
$ws.Range('E11').Value = '  +5.04%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('E20').Value = '  -0.29%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  -1.95%  '
$ws.Range('E27').Value = '  -4.12%  '
$ws.Range('E28').Value = '  -1.21%  '
$ws.Range('E29').Value = '  +0.44%  '
$ws.Range('E30').Value = '  +2.76%  '
$ws.Range('E31').Value = '  -0.94%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('E38').Value = '  -1.71%  '
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('E40').Value = '  -1.98%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('E42').Value = '  +0.85%  '
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('E44').Value = '  -0.62%  '
$ws.Range('E45').Value = '  +0.16%  '
$ws.Range('E46').Value = '  +3.48%  '
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('E50').Value = '  +3.40%  '
$ws.Range('E51').Value = '  -0.34%  '
